$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Restructure columns -------------------------------------------------
# Original layout: A=cst_code, B=cst_labe, C=updated_at
# Target layout:   A=_airbyte_ab_id, B=_airbyte_emitted_at, C=cst_code,
#                   D=cst_labe, E=_airbyte_additional_properties,
#                   F=source_file_path, G=updated_at

# Insert two new columns before the current column A
# (cst_code/cst_labe/updated_at shift from A:C to C:E)
$ws.Columns.Item(1).Insert()
$ws.Columns.Item(1).Insert()

# Insert two new columns right after cst_labe (current column D),
# before updated_at (which shifts from E to G)
$ws.Columns.Item(5).Insert()
$ws.Columns.Item(5).Insert()

# Give the brand-new header cells (A1,B1,E1,F1) the same formatting as the
# other header cells (bold / bordered / centered) by copying it over
$ws.Range("G1").Copy()
$ws.Range("A1:B1").PasteSpecial(-4122)
$ws.Range("G1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)

# --- Header row ------------------------------------------------------------
$ws.Cells.Item(1,1).Value = "_airbyte_ab_id"
$ws.Cells.Item(1,2).Value = "_airbyte_emitted_at"
$ws.Cells.Item(1,5).Value = "_airbyte_additional_properties"
$ws.Cells.Item(1,6).Value = "source_file_path"
# C1 (cst_code), D1 (cst_labe) and G1 (updated_at) already hold their
# correct text/style after the column inserts above.

# --- Data rows ---------------------------------------------------------
# (_airbyte_ab_id, source_file_path) per original data row, in row order
$rowData = @(
    @("160cca46-eac7-4a74-adb6-70d79ba5c70e"),
    @("ccb2d125-c888-4ab1-b409-66396db208b3"),
    @("3b51ce48-92d6-4cc2-81c9-3bb38aeb3cef"),
    @("4a6cadba-7c7b-43c4-b4d8-79202aca021a"),
    @("1771a560-e4f7-4236-a676-9bc3bf5ca205"),
    @("d6c73c51-daca-4145-aed0-c0ba09a9c6a1"),
    @("69cd510c-b6d3-4f1e-bd13-7401bab3fdfd"),
    @("c449afdf-100d-4ff3-a7c4-cbbbb4f091b4"),
    @("82c076eb-245c-4521-a5be-8bf4a75edae7"),
    @("03706e9c-5e75-4550-b5f5-122ee88ff115"),
    @("de28f25d-f037-4720-b46c-78e000ff5b58"),
    @("cb8569a9-fce2-40ff-afd0-36b2f72ea490"),
    @("e833be5e-0d48-44f5-a97b-3ea7be9740ed"),
    @("79481592-dbd4-4196-8a1a-a311b0c1666f"),
    @("95ac4d57-ab03-4936-aea5-f776b48bea64"),
    @("69c9ee75-366b-491b-9db7-1dac2e2365f7")
)

$sourceFile = "s3a://ai360nica/data/bronze/oracle/eftswitch/MXP/CARD_STATUS/2024_08_06_1722929004063_0.parquet"

for ($i = 0; $i -lt $rowData.Length; $i++) {
    $r = $i + 2
    $ws.Cells.Item($r,1).Value = $rowData[$i][0]

    $ws.Cells.Item($r,2).NumberFormat = "YYYY-MM-DD HH:MM:SS"
    $ws.Cells.Item($r,2).Value = 45510.3079196875

    $ws.Cells.Item($r,6).Value = $sourceFile

    $ws.Cells.Item($r,7).Value = 45511.29450409365
}
